# Commit: "update to consultancy mode"
#
# The "algorithms" sheet's used range shrank from 509 data rows down to 118
# (the _FilterDatabase now covers B1:G118 instead of B1:H509), and a new
# column H was populated with the marker value "d" for most of the
# remaining rows (51-118), except for a handful of rows that were
# deliberately left blank in column H.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("algorithms")
$ws.Activate()

# Rows (in the algorithms sheet) that get "d" written into column H.
$rowsWithD = @(
    51, 52, 53, 54, 55, 56, 57, 58, 59, 60, 61, 62, 63,
    66,
    68,
    70, 71, 72, 73, 74, 75, 76, 77, 78, 79, 80, 81,
    84,
    86,
    88, 89, 90,
    93, 94, 95, 96, 97, 98, 99, 100, 101, 102, 103, 104, 105, 106, 107, 108, 109, 110, 111, 112, 113, 114, 115, 116, 117, 118
)

foreach ($r in $rowsWithD) {
    $ws.Cells.Item($r, 8).Value = "d"
}

# Re-apply the autofilter over the now-shrunk data range (B1:G118) instead
# of the old B1:H509, and make sure it is reflected in the sheet's hidden
# _FilterDatabase defined name.
$ws.Range("B1:G118").AutoFilter() | Out-Null
$wb.Names.Item("algorithms!_FilterDatabase").RefersTo = "=algorithms!`$B`$1:`$G`$118"

# Update the cursor/selection position left in the saved view.
$ws.Range("L14").Select()

Write-Output "done"
